$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 gains a note about reading papers on feature selection
$ws.Range("D28").WrapText = $true
$ws.Range("D28").Value = "Reading papers on feature selection. Lasso, "

# New entry: Monday June 3rd, 2024 - regrouping with Daniel, starting a new
# rmarkdown file to work on GBMs
$ws.Range("A29").NumberFormat = "d-mmm"
$ws.Range("A29").Value = 45446
$ws.Range("B29").Value = 3
$ws.Range("D29").WrapText = $true
$ws.Range("D29").Value = "Regroup w Daniel. Starting new rmd on GBMs"

# Leave the selection on the newly added note, matching where work left off
$ws.Range("D29").Select()
